$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the tab name in A2 from "CasesTab" to "ParticipantsTab"
$ws.Range("A2").Value = "ParticipantsTab"

# Move the active selection from B4 to A2
$ws.Range("A2").Select()
